$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F3:F28").Borders(10).LineStyle = -4142
